# Revision Inicial Semana 1 con Carolina
# Incluye primera sesion de revision e incorporacion de comentarios a la
# semana 1 (revision con Carolina): reposition two pictures on slide 3.
#
# Target EMU offsets (size/ext unchanged):
#   Picture 10:  (8026137,3368655) -> (8138735,3024634)   ext 596433 x 596433
#   Picture 15:  (7559845,4241966) -> (7770484,3758707)   ext 831476 x 1110020
#
# Shape.Left/Top are expressed in points (1 pt = 12700 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$picture10 = $s.Shapes.Item("Picture 10")
$picture10.Left = 640.8452755905512
$picture10.Top = 238.16016388031494

$picture15 = $s.Shapes.Item("Picture 15")
$picture15.Left = 611.8491516582677
$picture15.Top = 295.9611811023622
